$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the two runs " when asking for " + "updates." into a single
#    run by doing a Find/Replace that spans the run boundary.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    " when asking for updates.", $false, $false, $false, $false, $false,
    $true, 1, $false, " when asking for updates.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Merge "The assignment is asking us to create unigram, bigram, and "
#    + "partial" + " trigram models. ... improvement." into one run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The assignment is asking us to create unigram, bigram, and partial trigram models. Additionally, we need to analyze the results obtained from these models and evaluate their performance on specific text. This involves examining how accurately the models predict the next word in a sequence and identifying any limitations or areas for improvement.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The assignment is asking us to create unigram, bigram, and partial trigram models. Additionally, we need to analyze the results obtained from these models and evaluate their performance on specific text. This involves examining how accurately the models predict the next word in a sequence and identifying any limitations or areas for improvement.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3) Merge "My main goal ... understanding " + "on NLP" +
#    ", enhance my skills, and collaborate effectively with my team
#    members" into one run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "My main goal for the project is to learn something new and gain valuable experience. While achieving a good grade is important, I see it as an outcome of the work rather than the primary goal. I aim to deepen my understanding on NLP, enhance my skills, and collaborate effectively with my team members",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "My main goal for the project is to learn something new and gain valuable experience. While achieving a good grade is important, I see it as an outcome of the work rather than the primary goal. I aim to deepen my understanding on NLP, enhance my skills, and collaborate effectively with my team members",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4) Add <w:lang w:val="en-US"/> to the paragraph mark run properties of
#    the paragraph holding "My personal development goal..." (Question 8
#    answer) by setting the paragraph mark's LanguageID.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "My personal development goal for this group project*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.LanguageID = 1033
}

# ---------------------------------------------------------------------
# 5) Append new paragraphs at the end of the document (before sectPr):
#      - blank paragraph
#      - "Question 9: What languages do you speak?"
#      - "I speak Arabic (mother tongue),English and a bit of Dutch "
#      - " " (single space)
# ---------------------------------------------------------------------
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd

$end.InsertParagraphAfter()
$end.Collapse(0)

$end.InsertAfter("Question 9: What languages do you speak?")
$end.InsertParagraphAfter()
$end.Collapse(0)

$end.InsertAfter("I speak Arabic (mother tongue),English and a bit of Dutch ")
$end.InsertParagraphAfter()
$end.Collapse(0)

$end.InsertAfter(" ")
